$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2025-12-29 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-30 Tuesday", 2)

$t = $d.Tables.Item(1)

# Row 1 (table row 1): 87÷9= 55÷7= 84÷3= 25÷8= 42÷5= -> 90÷3= 34÷7= 66÷4= 92÷4= 62÷8=
$r = $t.Rows.Item(1)
$r.Cells.Item(1).Range.Text = "90÷3="
$r.Cells.Item(2).Range.Text = "34÷7="
$r.Cells.Item(3).Range.Text = "66÷4="
$r.Cells.Item(4).Range.Text = "92÷4="
$r.Cells.Item(5).Range.Text = "62÷8="

# Row 2 (table row 5): 87÷9= 25÷4= 56÷2= 17÷9= 12÷7= -> 20÷2= 38÷6= 71÷3= 51÷6= 96÷5=
$r = $t.Rows.Item(5)
$r.Cells.Item(1).Range.Text = "20÷2="
$r.Cells.Item(2).Range.Text = "38÷6="
$r.Cells.Item(3).Range.Text = "71÷3="
$r.Cells.Item(4).Range.Text = "51÷6="
$r.Cells.Item(5).Range.Text = "96÷5="

# Row 3 (table row 9): 40÷4= 86÷8= 51÷3= 49÷3= 87÷6= -> 77÷4= 40÷4= 73÷8= 15÷8= 89÷8=
# (net effect of inserting a 77÷4= cell at the front and dropping the 49÷3= cell)
$r = $t.Rows.Item(9)
$r.Cells.Item(1).Range.Text = "77÷4="
$r.Cells.Item(2).Range.Text = "40÷4="
$r.Cells.Item(3).Range.Text = "73÷8="
$r.Cells.Item(4).Range.Text = "15÷8="
$r.Cells.Item(5).Range.Text = "89÷8="

# Row 4 (table row 13): 26÷2= 44÷3= 48÷8= 23÷3= 22÷4= -> 56÷4= 89÷7= 78÷6= 99÷5= 91÷4=
$r = $t.Rows.Item(13)
$r.Cells.Item(1).Range.Text = "56÷4="
$r.Cells.Item(2).Range.Text = "89÷7="
$r.Cells.Item(3).Range.Text = "78÷6="
$r.Cells.Item(4).Range.Text = "99÷5="
$r.Cells.Item(5).Range.Text = "91÷4="

# Row 5 (table row 17): 90÷9= 42÷8= 47÷8= 43÷7= 25÷7= -> 51÷3= 14÷2= 89÷2= 21÷7= 95÷8=
$r = $t.Rows.Item(17)
$r.Cells.Item(1).Range.Text = "51÷3="
$r.Cells.Item(2).Range.Text = "14÷2="
$r.Cells.Item(3).Range.Text = "89÷2="
$r.Cells.Item(4).Range.Text = "21÷7="
$r.Cells.Item(5).Range.Text = "95÷8="
